$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.315.75"
$ws.Range("E2").Value = "  +1.70%  "

$ws.Range("D3").Value = "1.895.13"
$ws.Range("E3").Value = "  -1.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.42"
$ws.Range("E5").Value = "  +1.64%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5158"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4023"
$ws.Range("E8").Value = "  +0.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08424"
$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("E10").Value = "  -0.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.115"
$ws.Range("E11").Value = "  -0.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.45"
$ws.Range("E12").Value = "  +11.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.439"
$ws.Range("E13").Value = "  +1.99%  "

$ws.Range("D14").Value = "1.896.14"
$ws.Range("E14").Value = "  -0.62%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.313"
$ws.Range("E15").Value = "  -0.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.63"
$ws.Range("E17").Value = "  +0.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001109"
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06646"
$ws.Range("E19").Value = "  -1.32%  "

$ws.Range("E20").Value = "  +1.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.941"
$ws.Range("E22").Value = "  -1.36%  "

$ws.Range("D23").Value = "30.294.63"
$ws.Range("E23").Value = "  +1.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.34"
$ws.Range("E24").Value = "  +1.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.225"
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("D26").Value = "2.119.28"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.75"
$ws.Range("E27").Value = "  +3.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.83"
$ws.Range("E28").Value = "  +2.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.320"
$ws.Range("E29").Value = "  -5.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.77"
$ws.Range("E30").Value = "  +0.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.084"
$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1050"
$ws.Range("E32").Value = "  -0.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.071"
$ws.Range("E33").Value = "  -1.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.844"
$ws.Range("E34").Value = "  +4.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02490"
$ws.Range("E35").Value = "  +0.13%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.356"
$ws.Range("E36").Value = "  +3.19%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06514"
$ws.Range("E37").Value = "  -1.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2199"
$ws.Range("E38").Value = "  -0.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.213"
$ws.Range("E39").Value = "  -2.26%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.89"
$ws.Range("E40").Value = "  +4.37%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.778"
$ws.Range("E41").Value = "  -2.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6482"
$ws.Range("E42").Value = "  -0.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.232"
$ws.Range("E43").Value = "  -0.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6081"
$ws.Range("E44").Value = "  -0.85%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.15"
$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.681"
$ws.Range("E46").Value = "  -0.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.057"
$ws.Range("E47").Value = "  -0.17%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.94"
$ws.Range("E48").Value = "  +0.19%  "

$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.233"
$ws.Range("E49").Value = "  -0.33%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "78.92"
$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.142"
$ws.Range("E51").Value = "  -2.02%  "
